$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "63.831.33"
$ws.Cells.Item(2,5).Value = "  -5.82%  "
$ws.Cells.Item(3,4).Value = "3.274.76"
$ws.Cells.Item(3,5).Value = "  -7.47%  "
$ws.Cells.Item(4,5).Value = "  -0.08%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "518.95"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "  -6.85%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "172.59"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "  -15.41%  "
$ws.Cells.Item(7,5).Value = "  -1.84%  "
$ws.Cells.Item(8,4).Value = "3.269.26"
$ws.Cells.Item(8,5).Value = "  -7.50%  "
$ws.Cells.Item(9,5).Value = "  +0.12%  "
$ws.Cells.Item(10,5).Value = "  -8.94%  "
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "55.95"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "  -13.29%  "
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.132"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "  -8.75%  "
$ws.Cells.Item(13,5).Value = "  -6.49%  "
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "8.95"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "  -10.24%  "
$ws.Cells.Item(15,4).Value = "3.788.22"
$ws.Cells.Item(15,5).Value = "  -7.25%  "
$ws.Cells.Item(16,4).Value = "3.266.75"
$ws.Cells.Item(16,5).Value = "  -7.33%  "
$ws.Cells.Item(17,5).Value = "  -7.14%  "
$ws.Cells.Item(18,4).Value = "63.659.60"
$ws.Cells.Item(18,5).Value = "  -5.71%  "
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "17.20"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "  -7.60%  "
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "10.94"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -8.13%  "
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "0.945"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "  -8.68%  "
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "369.96"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "  -6.36%  "
$ws.Cells.Item(23,5).Value = "  -7.78%  "
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "79.62"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = "  -4.68%  "
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "10.86"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "  -10.99%  "
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "3.84"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "  -2.10%  "
$ws.Cells.Item(27,2).Value = "LEO"
$ws.Cells.Item(27,3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "6.06"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "  -2.06%  "
$ws.Cells.Item(28,2).Value = "ImmutableX"
$ws.Cells.Item(28,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "2.63"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = "  -8.00%  "
$ws.Cells.Item(29,2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(29,3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "11.21"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = "  -9.15%  "
$ws.Cells.Item(30,2).Value = "Filecoin"
$ws.Cells.Item(30,3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "8.19"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = "  -8.40%  "
$ws.Cells.Item(31,2).Value = "EthereumClassic"
$ws.Cells.Item(31,3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "28.41"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = "  -9.11%  "
$ws.Cells.Item(32,2).Value = "Bittensor"
$ws.Cells.Item(32,3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "635.27"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = "  -11.67%  "
$ws.Cells.Item(33,2).Value = "NEARProtocol"
$ws.Cells.Item(33,3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "6.54"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = "  -8.83%  "
$ws.Cells.Item(34,2).Value = "Cosmos"
$ws.Cells.Item(34,3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "11.11"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "  -5.89%  "
$ws.Cells.Item(35,5).Value = "  -7.40%  "
$ws.Cells.Item(36,2).Value = "OKB"
$ws.Cells.Item(36,3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "58.51"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "  -8.96%  "
$ws.Cells.Item(37,2).Value = "Dai"
$ws.Cells.Item(37,3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "1.00"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = "  -0.15%  "
$ws.Cells.Item(38,2).Value = "InjectiveProtocol"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "36.19"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = "  -7.09%  "
$ws.Cells.Item(39,2).Value = "TheGraph"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.381"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = "  -5.23%  "
$ws.Cells.Item(40,2).Value = "FirstDigitalUSD"
$ws.Cells.Item(40,3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "0.996"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "  -0.14%  "
$ws.Cells.Item(41,2).Value = "PEPE"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(41,4).Value = "0.0₃0685"
$ws.Cells.Item(41,5).Value = "  -1.04%  "
$ws.Cells.Item(42,2).Value = "Maker"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(42,4).Value = "2.894.77"
$ws.Cells.Item(42,5).Value = "  -5.70%  "
$ws.Cells.Item(43,2).Value = "Kaspa"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.122"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "  -7.51%  "
$ws.Cells.Item(44,2).Value = "Fetch.AI"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "2.41"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "  -6.87%  "
$ws.Cells.Item(45,2).Value = "ThetaToken"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "2.64"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "  -13.02%  "
$ws.Cells.Item(46,2).Value = "WEMIXToken"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "2.62"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "  -4.78%  "
$ws.Cells.Item(47,2).Value = "VeChain"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "0.0391"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "  -4.97%  "
$ws.Cells.Item(48,2).Value = "ApeXProtocol"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "2.97"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "  +3.55%  "
$ws.Cells.Item(49,2).Value = "Stellar"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "0.124"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "  -3.16%  "
$ws.Cells.Item(50,2).Value = "Stacks"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "2.72"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "  +2.77%  "
$ws.Cells.Item(51,2).Value = "Monero"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "133.94"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "  -3.40%  "
